# rt_site_metadata.xlsx update:
#   - remove the two "roving receiver" rows (MT / RV) that used to sit at
#     rows 33-34, which pulls every pit_array row below them up by two
#   - drop the stray formatted-but-empty trailing row
#   - append two new pit_array sites for the Upper Salmon River (rkm 437 / 460)
#   - leave the selection on F47, matching the new last data cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the two roving-receiver rows (33: MT / Roving Receiver 17_18,
#    34: RV / Roving Receiver 19_20). Everything below shifts up by two,
#    so the old row 35 ("Lower Lemhi River") becomes row 33, etc.
$ws.Rows("33:34").Delete() | Out-Null

# 2) The old trailing row 48 (just a formatted, empty D48 cell) is now row
#    46 after the shift above - delete it so we can lay down real data.
$ws.Rows("46:46").Delete() | Out-Null

# 3) Seed rows 46-47 with the same formatting as the row above (row 45)
#    before writing the new site records into them.
$ws.Range("A45:J45").Copy() | Out-Null
$ws.Range("A46:J46").PasteSpecial(-4122) | Out-Null
$ws.Range("A47:J47").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 4) New site: Upper Salmon River at rkm 437 (USE)
$ws.Cells.Item(46, 1).Value2 = "Upper Salmon River at rkm 437"
$ws.Cells.Item(46, 2).Value2 = "USE"
$ws.Cells.Item(46, 3).Value2 = "pit_array"
$ws.Cells.Item(46, 4).Value2 = "NA"
$ws.Cells.Item(46, 5).Value2 = 45.028530000000003
$ws.Cells.Item(46, 6).Value2 = -113.916319
$ws.Cells.Item(46, 7).Value2 = $true
$ws.Cells.Item(46, 8).Value2 = $true
$ws.Cells.Item(46, 9).Value2 = $true
$ws.Cells.Item(46, 10).Value2 = "522.303.437"

# 5) New site: Upper Salmon River at rkm 460 (USI)
$ws.Cells.Item(47, 1).Value2 = "Upper Salmon River at rkm 460"
$ws.Cells.Item(47, 2).Value2 = "USI"
$ws.Cells.Item(47, 3).Value2 = "pit_array"
$ws.Cells.Item(47, 4).Value2 = "NA"
$ws.Cells.Item(47, 5).Value2 = 44.889763000000002
$ws.Cells.Item(47, 6).Value2 = -113.964145
$ws.Cells.Item(47, 7).Value2 = $true
$ws.Cells.Item(47, 8).Value2 = $true
$ws.Cells.Item(47, 9).Value2 = $true
$ws.Cells.Item(47, 10).Value2 = "522.303.460"

# 6) Match the saved selection from the authored workbook.
$ws.Range("F47").Select() | Out-Null
